$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.295.29"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.664.09"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "'218.91"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'0.5352"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "'0.2657"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'0.06405"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'0.07831"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'4.565"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "1.676.15"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "1.891.94"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "'0.5526"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "0.0₅8206"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "'65.75"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'4.679"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "'193.98"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'146.15"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "'0.1232"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "'7.193"
$ws.Range("D27").Value = "'16.10"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'1.485"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "'0.05848"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "'1.286"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "'3.605"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'3.283"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "'0.9642"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").Value = "'0.01606"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'0.8670"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").Value = "'5.875"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").Value = "1.052.99"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'104.75"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "1.802.51"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "'57.90"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("D47").Value = "'1.014"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "'0.4385"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "'8.009"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "'0.05164"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'1.414"
$ws.Range("E51").Value = "  -4.07%  "
